$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set string values in the exact order needed to reproduce shared string table order:
# 18: Items/Icons/Consumable/Potion_1 (D2)
# 19: Power Potion (B3)
# 20: 공격력 증가 (C3)
# 21: Items/Icons/Consumable/Potion_2 (D3)
# 22: Normal (I2)
$ws.Range("D2").Value = "Items/Icons/Consumable/Potion_1"
$ws.Range("B3").Value = "Power Potion"
$ws.Range("C3").Value = "공격력 증가"
$ws.Range("D3").Value = "Items/Icons/Consumable/Potion_2"
$ws.Range("I2").Value = "Normal"
$ws.Range("I3").Value = "Normal"

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = 20
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0

$ws.Range("A3").Value = 10301021
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("H3").Value = 20
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0

$ws.Columns.Item(4).ColumnWidth = 30.428571428571427
$ws.Columns.Item(7).ColumnWidth = 12

$ws.Range("K14").Select()

$wb.Save()
